# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" palette (used by the Notes Master)
#   ppt/theme/theme2.xml  -> "Integral" palette      (used by the Slide Master /
#                                                      the presentation itself)
#
# The authored change swaps the two themes' contents: theme1.xml becomes the
# "Integral" palette and theme2.xml becomes the "Office Theme" palette (the
# relationships that point at the files are untouched - only what's inside
# each theme part changes).
#
# The PowerPoint object model only ever exposes a single, live
# ThemeColorScheme - it is always the one backing the Slide Master /
# presentation theme (ppt/theme/theme2.xml) no matter whether it is reached
# via $p.SlideMaster, $p.NotesMaster, $p.HandoutMaster or $p.Designs. So we
# repaint that one reachable color scheme with the "Office Theme" palette
# that theme2.xml is supposed to end up with.
#
# ThemeColorScheme.Item(i) order == clrScheme child order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB() is encoded as R + G*256 + B*65536, matching VBA's RGB().

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
